$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the Heading1 title, re-using its bold run to build a
#    brand-new "page title" paragraph placed just above the final
#    "Prompt: ..." paragraph at the end of the document.
# ------------------------------------------------------------------

$metaPara = $d.Paragraphs(2)
$metaPara.Range.Cut()

$lastIdx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIdx)
$insertion = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertion.Paste()

# The pasted paragraph now occupies the slot the last paragraph used
# to be at; the "Prompt" paragraph got pushed one slot further down.
$pastedIdx = $lastIdx
$pastedPara = $d.Paragraphs($pastedIdx)

# Split the pasted paragraph right after the bold "Meta description"
# run so the label and the description text land in separate
# paragraphs.
$splitRange = $pastedPara.Range.Duplicate()
$splitRange.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($splitRange.End, $splitRange.End)
$splitPoint.InsertParagraphAfter()

$titlePara = $d.Paragraphs($pastedIdx)
$descPara = $d.Paragraphs($pastedIdx + 1)

# Re-purpose the bold run's text into the page-title paragraph that
# now precedes the (still present, untouched) "Prompt" paragraph.
$titleRange = $titlePara.Range.Duplicate()
$titleRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "Play 5 Lucky Lions Free Slot Game | Review & RTP 96.79%", 2)

# The leftover ": Experience ..." paragraph is no longer needed as a
# standalone paragraph -- its text is moved into the existing "Prompt"
# paragraph below -- so delete it.
$descPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Replace the old "Prompt: ..." text with the (colon-stripped)
#    description text, keeping that paragraph's own formatting/markup
#    completely untouched otherwise.
# ------------------------------------------------------------------
$promptIdx = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($promptIdx)
$promptRange = $promptPara.Range.Duplicate()
$promptRange.Find.Execute("Prompt: Create a fun and vibrant feature image for ""5 Lucky Lions"" that captures the spirit of Chinese New Year and the colorful world of online slots. The image should feature a happy Maya warrior with glasses, as well as the game's symbols, including the lions, bonze, and drum. The overall style should be cartoonish and lively, evoking the festive atmosphere of lion dances and traditional celebrations. The image should be colorful and eye-catching, with a focus on the game's key features, such as the reels, paylines, and bonus pick feature.", $true, $false, $false, $false, $false, $true, 1, $false, "Experience the Chinese tradition of lion dance playing 5 Lucky Lions free slot game with 88 paylines and RTP 96.79%. Review and bets as low as €0.01.", 2)

Write-Output "ok"
